# Add two new test-case rows to the "Test Cases" sheet (Suite C), per the
# "Added new tests to Suite C" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 69 (E69) already carries the plain bordered style (s="1") that the new
# rows should use for every column, so copy its format across the new cells
# before writing values into them.
$ws.Range("E69").Copy()
$ws.Range("A70:E71").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New test case: AddExternalLinksToComments
$ws.Range("A70").Value = "AddExternalLinksToComments"
$ws.Range("B70").Value = "OPQA-1092"
$ws.Range("C70").Value = "Verfiy that user is able to add external links to the comment"
$ws.Range("D70").Value = "Y"
$ws.Range("E70").Value = "PASS"

# New test case: AddInternalLinksToComments
$ws.Range("A71").Value = "AddInternalLinksToComments"
$ws.Range("B71").Value = "OPQA-1093"
$ws.Range("C71").Value = "Verfiy that user is able to add links other NEON content [ex -Posts, articles, patents, profiles] to the comment"
$ws.Range("D71").Value = "Y"
$ws.Range("E71").Value = "PASS"

# Move the selection to reflect where the editor ended up after adding rows.
$ws.Range("B73").Select() | Out-Null
